$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.590.69'
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = '2.459.30'
$ws.Range("E3").Value = '  -1.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.91'
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.42'
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("E7").Value = '  +1.73%  '

$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("E9").Value = '  +3.63%  '

$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0796'
$ws.Range("E11").Value = '  +2.80%  '

$ws.Range("E12").Value = '  +0.80%  '

$ws.Range("D13").Value = '2.837.50'
$ws.Range("E13").Value = '  -0.96%  '

$ws.Range("E14").Value = '  +0.93%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.84'
$ws.Range("E15").Value = '  +3.62%  '

$ws.Range("D16").Value = '2.459.80'
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.775'
$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("D18").Value = '41.618.13'
$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.46'
$ws.Range("E19").Value = '  +2.80%  '

$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  +2.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.76'
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.35'
$ws.Range("E22").Value = '  +2.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.05'
$ws.Range("E23").Value = '  +1.61%  '

$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  +1.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.29'
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  +0.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.04'
$ws.Range("E30").Value = '  -3.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.97'
$ws.Range("E31").Value = '  +1.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.46'
$ws.Range("E32").Value = '  +1.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.57'
$ws.Range("E33").Value = '  +0.49%  '

$ws.Range("E34").Value = '  +0.24%  '

$ws.Range("B35").Value = 'ApeXProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.46'
$ws.Range("E35").Value = '  -0.71%  '

$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.42'
$ws.Range("E36").Value = '  -3.77%  '

$ws.Range("E37").Value = '  -3.31%  '

$ws.Range("E38").Value = '  +1.66%  '

$ws.Range("E39").Value = '  +1.92%  '

$ws.Range("E40").Value = '  -2.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.95'
$ws.Range("E41").Value = '  -2.56%  '

$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("D43").Value = '1.972.74'
$ws.Range("E43").Value = '  +1.27%  '

$ws.Range("E44").Value = '  +0.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.70'
$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.91'
$ws.Range("E46").Value = '  -1.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.91'
$ws.Range("E47").Value = '  +2.09%  '

$ws.Range("D48").Value = '2.695.08'
$ws.Range("E48").Value = '  -0.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.42'
$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.63'
$ws.Range("E50").Value = '  +0.24%  '

$ws.Range("E51").Value = '  -1.95%  '
